# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh
# (GitHub Actions symbol-list update, Mon Feb  6 19:24:39 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '329.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.30%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.20%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.592'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '3.40%'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.09%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.036'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '6.61%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.303'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.04%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9536'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.32%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-6.03%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1185'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.08%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1854'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.43%'
$ws.Range("B12").Value = 'MCDex'
$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '10.20'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '20.13%'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09729'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.16%'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04595'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '10.26%'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1069'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.13%'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001274'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.10%'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04214'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-3.80%'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005852'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-4.02%'
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.369'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-5.63%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3474'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.71%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.60%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2502'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.93%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001247'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.70%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004325'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.21%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001188'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.90%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.80%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02675'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-0.04%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05567'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.18%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007574'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.94%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1411'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.26%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.008068'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-17.36%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002013'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.18%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008396'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-15.71%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007205'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.93%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.67%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.004715'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '35.97%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002267'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002097'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.67%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001997'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.67%'
